# "Added column in QC Text.xlsx file"
#
# The Controller sheet gets a new "DropDown Select2" column inserted right
# before the existing "Delete" column, so "Delete" (and its data) shifts
# one column to the right (K -> L). The new column carries a header in
# row 1 and a "No" value in row 2, with row 3 left blank - formatted like
# their neighbouring cells in column J.

$wb = $excel.ActiveWorkbook
$wsModel = $wb.Worksheets.Item("Model")
$wsViews = $wb.Worksheets.Item("Views")
$wsCtrl = $wb.Worksheets.Item("Controller")

# --- Controller sheet: insert the new column ---------------------------
$wsCtrl.Columns.Item(11).Insert()

$wsCtrl.Cells.Item(1, 11).Value = "DropDown Select2"
$wsCtrl.Cells.Item(1, 11).Font.Bold = $wsCtrl.Cells.Item(1, 10).Font.Bold
$wsCtrl.Cells.Item(1, 11).Borders.LineStyle = $wsCtrl.Cells.Item(1, 10).Borders.LineStyle

$wsCtrl.Cells.Item(2, 11).Value = "No"
$wsCtrl.Cells.Item(2, 11).Font.Bold = $wsCtrl.Cells.Item(2, 10).Font.Bold
$wsCtrl.Cells.Item(2, 11).Borders.LineStyle = $wsCtrl.Cells.Item(2, 10).Borders.LineStyle

$wsCtrl.Cells.Item(3, 11).Font.Bold = $wsCtrl.Cells.Item(3, 10).Font.Bold
$wsCtrl.Cells.Item(3, 11).Borders.LineStyle = $wsCtrl.Cells.Item(3, 10).Borders.LineStyle

# Widen the new column to fit its header text (matches the bestFit width
# Excel would compute for "DropDown Select2").
$wsCtrl.Columns.Item(11).ColumnWidth = 16.71

# --- Selections / active sheet -----------------------------------------
$wsModel.Range("H2").Select()
$wsViews.Range("E2").Select()

$wsCtrl.Activate()
$wsCtrl.Range("H11").Select()
